$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Done" -> "DONE" rows (simple rename, status unchanged semantically)
$doneRows = @(2, 3, 4, 19, 20, 31)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 5).Value = "DONE"
}

# Rows that moved from "Made" + F "DUE" to fully "DONE" (F cleared)
$completedRows = @(13, 21, 22, 23, 26, 32)
foreach ($r in $completedRows) {
    $ws.Cells.Item($r, 5).Value = "DONE"
    $ws.Cells.Item($r, 6).Value = ""
}

# Rows that gained a "DUE" marker in column F (previously blank)
$newDueRows = @(11, 12, 14)
foreach ($r in $newDueRows) {
    $ws.Cells.Item($r, 6).Value = "DUE"
}

# Rows that lost their "DUE" marker in column F (now blank)
$removedDueRows = @(29, 30, 33)
foreach ($r in $removedDueRows) {
    $ws.Cells.Item($r, 6).Value = ""
}

# Row 28: status becomes DONE while F (DUE) stays
$ws.Cells.Item(28, 5).Value = "DONE"

# Update the selected cell in the sheet view
$ws.Range("F21").Select()
